# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets for rows 2-7.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row number -> new value for column F
$updates = @{
    2 = 2178
    3 = 624
    4 = 1558
    5 = 7316
    6 = 179
    7 = 174
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
